# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to reflect newly generated site output.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 13901
$wsExpo.Range("F6").Value = 496
$wsExpo.Range("F10").Value = 14849
$wsExpo.Range("F26").Value = 5769

# --- Sheet: 全部类型 (all types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 13901
$wsAll.Range("F7").Value = 496
$wsAll.Range("F11").Value = 14849
$wsAll.Range("F27").Value = 5769
